$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '25.815.20'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -0.02%  '
$ws.Range('E2').Style = "Normal"

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.739.34'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -1.00%  '
$ws.Range('E3').Style = "Normal"

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('E4').Style = "Normal"

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '225.53'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -4.80%  '
$ws.Range('E5').Style = "Normal"

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.001'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.11%  '
$ws.Range('E6').Style = "Normal"

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.5182'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +2.56%  '
$ws.Range('E7').Style = "Normal"

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2730'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  +2.88%  '
$ws.Range('E8').Style = "Normal"

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '38.96'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -5.38%  '
$ws.Range('E9').Style = "Normal"

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.06071'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -2.04%  '
$ws.Range('E10').Style = "Normal"

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.743.47'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -0.42%  '
$ws.Range('E11').Style = "Normal"

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.07003'
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +1.01%  '
$ws.Range('E12').Style = "Normal"

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '15.12'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -3.34%  '
$ws.Range('E13').Style = "Normal"

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.6306'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +5.01%  '
$ws.Range('E14').Style = "Normal"

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.481'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.19%  '
$ws.Range('E15').Style = "Normal"

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '76.27'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -1.50%  '
$ws.Range('E16').Style = "Normal"

$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +0.18%  '
$ws.Range('E17').Style = "Normal"

$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +0.05%  '
$ws.Range('E18').Style = "Normal"

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '25.835.11'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -0.05%  '
$ws.Range('E19').Style = "Normal"

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '11.43'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -1.51%  '
$ws.Range('E20').Style = "Normal"

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '0.000006588'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -3.24%  '
$ws.Range('E21').Style = "Normal"

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.961.03'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -0.71%  '
$ws.Range('E22').Style = "Normal"

$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +0.03%  '
$ws.Range('E23').Style = "Normal"

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '8.398'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +2.18%  '
$ws.Range('E24').Style = "Normal"

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '5.081'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  -2.22%  '
$ws.Range('E25').Style = "Normal"

$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '136.72'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -0.81%  '
$ws.Range('E26').Style = "Normal"

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.491'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +2.60%  '
$ws.Range('E27').Style = "Normal"

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.811'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -0.73%  '
$ws.Range('E28').Style = "Normal"

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '14.93'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -0.46%  '
$ws.Range('E29').Style = "Normal"

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '102.58'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +0.06%  '
$ws.Range('E30').Style = "Normal"

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.08287'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +1.31%  '
$ws.Range('E31').Style = "Normal"

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.610'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -1.58%  '
$ws.Range('E32').Style = "Normal"

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.363'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -0.50%  '
$ws.Range('E33').Style = "Normal"

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.04389'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -0.20%  '
$ws.Range('E34').Style = "Normal"

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '2.626'
$ws.Range('D35').Style = "Normal"

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.9634'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -3.49%  '
$ws.Range('E36').Style = "Normal"

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.5943'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -1.15%  '
$ws.Range('E37').Style = "Normal"

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.676'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -1.51%  '
$ws.Range('E38').Style = "Normal"

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01553'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +0.29%  '
$ws.Range('E39').Style = "Normal"

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.924'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -0.57%  '
$ws.Range('E40').Style = "Normal"

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.000'
$ws.Range('D41').Style = "Normal"

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '100.90'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -2.15%  '
$ws.Range('E42').Style = "Normal"

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.3798'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -0.15%  '
$ws.Range('E43').Style = "Normal"

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.7281'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -0.81%  '
$ws.Range('E44').Style = "Normal"

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '4.846'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -1.67%  '
$ws.Range('E45').Style = "Normal"

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.05482'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -0.14%  '
$ws.Range('E46').Style = "Normal"

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '6.153'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  +3.80%  '
$ws.Range('E47').Style = "Normal"

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.1098'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +0.80%  '
$ws.Range('E48').Style = "Normal"

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '29.64'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -0.37%  '
$ws.Range('E49').Style = "Normal"

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '51.87'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -0.19%  '
$ws.Range('E50').Style = "Normal"

$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  +0.26%  '
$ws.Range('E51').Style = "Normal"
